$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.841.12"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.083.75"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.07"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.89"
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "2.391.40"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.76"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.22"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.768"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "2.082.22"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "37.753.21"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.31"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.02"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.45"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +7.69%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.55"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.41"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.70"
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0992"
$ws.Range("E40").Value = "  +4.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.84"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.44"
$ws.Range("E43").Value = "  +8.13%  "
$ws.Range("D44").Value = "1.467.61"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0215"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.05"
$ws.Range("E47").Value = "  +5.50%  "
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").Value = "2.276.16"
$ws.Range("E51").Value = "  +0.55%  "
